$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: new datetime + new price
$ws.Range("C2").Value = 44180.38637571139
$ws.Range("E2").Value = 100

# Update row 3: new datetime + new item name + new price
$ws.Range("C3").Value = 44180.38637571139
$ws.Range("D3").Value = "подстаканник"
$ws.Range("E3").Value = 200

# Remove rows 4-9 entirely (data shrinks to A1:E3)
$ws.Range("A4:E9").EntireRow.Delete()
